$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename Sheet2 -> SHTeam2016 ---
$ws2.Name = "SHTeam2016"

# --- Copy Sheet1's header+data block into SHTeam2016, row by row,
#     formats first then values, so the engine dedups against the
#     existing style table (s=1..8) instead of minting new styles. ---
for ($r = 1; $r -le 6; $r++) {
    $srcRow = $ws1.Range("A" + $r + ":G" + $r)
    $dstRow = $ws2.Range("A" + $r + ":G" + $r)
    $srcRow.Copy()
    $dstRow.PasteSpecial(-4122)   # xlPasteFormats
}
for ($r = 1; $r -le 6; $r++) {
    $srcRow = $ws1.Range("A" + $r + ":G" + $r)
    $dstRow = $ws2.Range("A" + $r + ":G" + $r)
    $srcRow.Copy()
    $dstRow.PasteSpecial(-4163)   # xlPasteValues
}

# --- Row heights (custom row heights carried over from Sheet1) ---
$ws2.Rows.Item(1).RowHeight = 14.25
$ws2.Rows.Item(2).RowHeight = 147.75
$ws2.Rows.Item(3).RowHeight = 74.25
$ws2.Rows.Item(4).RowHeight = 255
$ws2.Rows.Item(5).RowHeight = 81
$ws2.Rows.Item(6).RowHeight = 66

# --- Column widths (carried over from Sheet1) ---
$ws2.Columns.Item(1).ColumnWidth = $ws1.Columns.Item(1).ColumnWidth
$ws2.Columns.Item(2).ColumnWidth = $ws1.Columns.Item(2).ColumnWidth
$ws2.Columns.Item(3).ColumnWidth = $ws1.Columns.Item(3).ColumnWidth
$ws2.Columns.Item(4).ColumnWidth = $ws1.Columns.Item(4).ColumnWidth
$ws2.Columns.Item(5).ColumnWidth = $ws1.Columns.Item(5).ColumnWidth
$ws2.Columns.Item(6).ColumnWidth = $ws1.Columns.Item(6).ColumnWidth
$ws2.Columns.Item(7).ColumnWidth = $ws1.Columns.Item(7).ColumnWidth

# --- Selections ---
# Sheet1: selection moves to E10, but the whole sheet is highlighted.
$ws1.Cells.Select()
# SHTeam2016: single-cell selection on D4, and this sheet becomes the active tab.
$ws2.Range("D4").Select()
$ws2.Activate()
